$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$sh = $m.Shapes.Item(4)
$tr2 = $sh.TextFrame2.TextRange
$paras = $tr2.Paragraphs()
Write-Host "paras count: $($paras.Count)"
for ($i=1; $i -le $paras.Count; $i++) {
  $para = $paras.Item($i)
  Write-Host "para $i text: '$($para.Text)'"
}
